# Trade #69 closed at 2026-02-17 15:48:12 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade (#69 / row 70).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.16   # Current Capital
$summary.Range("B4").Value = 0.15      # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 69        # Total Trades
$summary.Range("B8").Value = 37        # Losing Trades
$summary.Range("B9").Value = 30.43     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.16     # Capital
$status.Range("D4").Value = 69         # Trades
$status.Range("E4").Value = 0.15       # P&L $
$status.Range("F4").Value = 0.16       # P&L %
$status.Range("G4").Value = 30.43      # Win Rate %

# ---------------------------------------------------------------------
# All Trades + MarketMaking sheets: append the newly closed trade as row 70
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A70").Value = 69

    # The Date column holds plain text like "2026-02-17"; format the cell as
    # text first so Excel doesn't auto-convert the literal into a date serial.
    $ws.Range("B70").NumberFormat = "@"
    $ws.Range("B70").Value = "2026-02-17"

    $ws.Range("C70").Value = "15:48:06"
    $ws.Range("D70").Value = "MarketMaking"
    $ws.Range("E70").Value = "DOWN"
    $ws.Range("F70").Value = 0.38
    $ws.Range("G70").Value = 0.373617
    $ws.Range("H70").Value = "CLOSED"
    $ws.Range("I70").Value = -1.6798
    $ws.Range("J70").Value = -0.01
    $ws.Range("K70").Value = 100.16
    $ws.Range("L70").Value = 0
    $ws.Range("M70").Value = 0
    $ws.Range("N70").Value = 0.6
    $ws.Range("O70").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P70").Value = "early_exit"
    $ws.Range("Q70").Value = 0.13
}
